# Applies the changes described in the commit:
# "Caricata prima versione senza classi e senza lettura file pdf"
#
# - Adds a new "docente" label in B5
# - Lower-cases the header labels "Progetto uno" / "Progetto due" in C1/D1
#   (now "progetto uno" / "progetto due")
# - Moves the active selection to D1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "docente" cell first so it gets the lowest free shared-string
# index, matching the order new strings were introduced upstream.
$ws.Range("B5").Value = "docente"

# Update the existing project headers to lowercase wording.
$ws.Range("C1").Value = "progetto uno"
$ws.Range("D1").Value = "progetto due"

# Move the selection/active cell to D1.
$ws.Range("D1").Select()
